# Insert a new data row at row 103 (pushing existing rows 103-161 down to
# 104-162) and populate it with the new Zapallo italiano price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(103).Insert()

$ws.Cells.Item(103, 1).Value2  = 11
$ws.Cells.Item(103, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(103, 3).Value2  = "Bíobío"
$ws.Cells.Item(103, 4).Value2  = 44873
$ws.Cells.Item(103, 5).Value2  = 8
$ws.Cells.Item(103, 6).Value2  = 100112032
$ws.Cells.Item(103, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(103, 8).Value2  = "Sin especificar"
$ws.Cells.Item(103, 9).Value2  = "Primera"
$ws.Cells.Item(103, 10).Value2 = 450
$ws.Cells.Item(103, 11).Value2 = 5500
$ws.Cells.Item(103, 12).Value2 = 6000
$ws.Cells.Item(103, 13).Value2 = 5722
$ws.Cells.Item(103, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(103, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(103, 16).Value2 = 114
$ws.Cells.Item(103, 17).Value2 = 50
$ws.Cells.Item(103, 18).Value2 = "Hortaliza"
